$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 275, shifting rows 275:378 down to 276:379.
$ws.Rows.Item(275).Insert()

# Populate the new row 275 with the new weekly price-report entry.
$ws.Range("A275").Value = 4
$ws.Range("B275").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C275").Value = "Los Lagos"
$ws.Range("D275").Value = 45093
$ws.Range("E275").Value = 10
$ws.Range("F275").Value = "Fruta"
$ws.Range("G275").Value = 100109
$ws.Range("H275").Value = "Uva"
$ws.Range("I275").Value = 100109001
$ws.Range("J275").Value = "Uva"
$ws.Range("K275").Value = "Red Globe"
$ws.Range("L275").Value = "Primera"
$ws.Range("M275").Value = 400
$ws.Range("N275").Value = 15000
$ws.Range("O275").Value = 15500
$ws.Range("P275").Value = 15250
$ws.Range("Q275").Value = "$/caja 18 kilos"
$ws.Range("R275").Value = "Región de O'Higgins"
$ws.Range("S275").Value = 847
$ws.Range("T275").Value = 18
